# Adds a new "2022-Q4" quarterly sheet right after the "总计" (summary) sheet,
# populates it with that quarter's fund-holdings data, and inserts a matching
# summary row at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write one fund-holdings data row.
#   col A -> numeric row index (0-based)
#   col B -> fund code            (text, keeps leading zeros)
#   col C -> fund name             (text)
#   col D -> fund size             (text, e.g. "26.05")
#   col E -> total stock position  (text, e.g. "91.85")
#   col F -> position ratio        (text, e.g. "5.92")
#   col G -> held market value     (text, e.g. "1.5422")
#   col H -> position rank         (number)
# ---------------------------------------------------------------------------
function Set-FundRow {
    param($ws, $row, $idxA, $code, $name, $size, $pos, $ratio, $mval, $rank)

    $ws.Cells.Item($row, 1).Value = $idxA

    $c = $ws.Cells.Item($row, 2); $c.Value = "'" + $code
    $c = $ws.Cells.Item($row, 3); $c.Value = "'" + $name
    $c = $ws.Cells.Item($row, 4); $c.Value = "'" + $size
    $c = $ws.Cells.Item($row, 5); $c.Value = "'" + $pos
    $c = $ws.Cells.Item($row, 6); $c.Value = "'" + $ratio
    $c = $ws.Cells.Item($row, 7); $c.Value = "'" + $mval

    $ws.Cells.Item($row, 8).Value = $rank
}

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet right after "总计" and name it "2022-Q4".
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Add($null, $wsTotal)
$ws4.Name = "2022-Q4"

# Header row (bold, centered, top-aligned, thin border) - same look as the
# other quarterly sheets.
$header = $ws4.Range("B1:H1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1

$ws4.Range("B1").Value = "基金代码"
$ws4.Range("C1").Value = "基金名称"
$ws4.Range("D1").Value = "基金规模"
$ws4.Range("E1").Value = "股票总仓位"
$ws4.Range("F1").Value = "仓位占比"
$ws4.Range("G1").Value = "持有市值(亿元)"
$ws4.Range("H1").Value = "仓位排名"

# Column A header style matches the rest of column A in this sheet.
$aHeader = $ws4.Range("A1")
$aHeader.Font.Bold = $true
$aHeader.HorizontalAlignment = -4108
$aHeader.VerticalAlignment = -4160
$aHeader.Borders.LineStyle = 1

# Data rows (fund holdings for 2022-Q4).
Set-FundRow $ws4 2 0 "000362" "国泰聚信价值优势灵活配置混合A" "26.05" "91.85" "5.92" "1.5422" 3
Set-FundRow $ws4 3 1 "020010" "国泰金牛创新混合" "13.29" "89.45" "5.87" "0.7801" 4
Set-FundRow $ws4 4 2 "000363" "国泰聚信价值优势灵活配置混合C" "12.90" "91.85" "5.92" "0.7637" 3
Set-FundRow $ws4 5 3 "012173" "国泰兴泽优选一年持有期混合A" "8.18" "92.65" "4.54" "0.3714" 6
Set-FundRow $ws4 6 4 "008415" "国泰大制造两年持有期混合" "9.64" "92.78" "3.84" "0.3702" 10
Set-FundRow $ws4 7 5 "011335" "银河医药健康混合A" "8.18" "92.65" "3.98" "0.3256" 9
Set-FundRow $ws4 8 6 "011466" "兴业医疗保健混合A" "4.60" "88.16" "7.05" "0.3243" 2
Set-FundRow $ws4 9 7 "012160" "财通资管健康产业混合C" "4.86" "94.25" "5.93" "0.2882" 7
Set-FundRow $ws4 10 8 "012174" "国泰兴泽优选一年持有期混合C" "5.86" "92.65" "4.54" "0.2660" 6
Set-FundRow $ws4 11 9 "013890" "国泰睿毅三年持有期混合A" "4.82" "90.04" "4.89" "0.2357" 4
Set-FundRow $ws4 12 10 "012159" "财通资管健康产业混合A" "3.89" "94.25" "5.93" "0.2307" 7
Set-FundRow $ws4 13 11 "002938" "中银证券健康产业灵活配置混合" "5.33" "92.95" "4.25" "0.2265" 4
Set-FundRow $ws4 14 12 "005244" "国泰聚优价值灵活配置混合A" "4.39" "90.86" "3.65" "0.1602" 5
Set-FundRow $ws4 15 13 "011467" "兴业医疗保健混合C" "2.12" "88.16" "7.05" "0.1495" 2
Set-FundRow $ws4 16 14 "008619" "永赢医药健康股票C" "1.78" "89.12" "7.67" "0.1365" 5
Set-FundRow $ws4 17 15 "005245" "国泰聚优价值灵活配置混合C" "2.10" "90.86" "3.65" "0.0766" 5
Set-FundRow $ws4 18 16 "151002" "银河收益混合" "9.01" "23.82" "0.73" "0.0658" 8
Set-FundRow $ws4 19 17 "008618" "永赢医药健康股票A" "0.58" "89.12" "7.67" "0.0445" 5
Set-FundRow $ws4 20 18 "013891" "国泰睿毅三年持有期混合C" "0.45" "90.04" "4.89" "0.0220" 4
Set-FundRow $ws4 21 19 "004914" "中银证券聚瑞混合C" "0.31" "42.02" "1.70" "0.0053" 10
Set-FundRow $ws4 22 20 "004913" "中银证券聚瑞混合A" "0.20" "42.02" "1.70" "0.0034" 10
Set-FundRow $ws4 23 21 "015666" "银河医药健康混合C" "0.02" "92.65" "3.98" "0.0008" 9

# ---------------------------------------------------------------------------
# 2. Rewrite the "总计" summary table: a new 2022-Q4 row goes on top (row 2)
#    and every other quarter's row shifts down by one. Row 10 is brand new,
#    so its column-A cell needs the same bold/centered/bordered styling that
#    the other column-A cells already carry.
# ---------------------------------------------------------------------------
$summaryRows = @(
    @(0, "2022-Q4", 22, 6.39),
    @(1, "2022-Q3", 4, 0.43),
    @(2, "2022-Q2", 9, 0.83),
    @(3, "2022-Q1", 7, 1.9),
    @(4, "2021-Q4", 5, 1),
    @(5, "2021-Q3", 1, 0.79),
    @(6, "2021-Q2", 6, 1.89),
    @(7, "2021-Q1", 7, 4.15),
    @(8, "2020-Q4", 4, 0.15)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $row = $i + 2
    $data = $summaryRows[$i]
    $wsTotal.Cells.Item($row, 1).Value = $data[0]
    $wsTotal.Cells.Item($row, 2).Value = $data[1]
    $wsTotal.Cells.Item($row, 3).Value = $data[2]
    $wsTotal.Cells.Item($row, 4).Value = $data[3]
}

# Row 10 did not exist before this edit - give its column-A cell the same
# bold/centered/top-aligned/thin-bordered look used by the rest of column A.
$aNew = $wsTotal.Cells.Item(10, 1)
$aNew.Font.Bold = $true
$aNew.HorizontalAlignment = -4108   # xlCenter
$aNew.VerticalAlignment = -4160     # xlTop
$aNew.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3. Keep "2020-Q4" (the last tab) as the selected/active sheet, same as
#    before the edit - adding the new sheet must not steal the selection.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
